$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Img size column (G) for existing rows 2 and 3: numeric 300 -> text "RandomResizedCrop(300)" ---
$ws.Range("G2").Value = "RandomResizedCrop(300)"
$ws.Range("G3").Value = "RandomResizedCrop(300)"

# --- Copy formatting of row 3 into new row 4 ---
$ws.Range("A3:O3").Copy()
$ws.Range("A4:O4").PasteSpecial(-4122)

# --- Renumber Sno column ---
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# --- Fill new row 4 data (Resnet18_imagenet_full model) ---
$ws.Range("B4").Value = "R18_imagenet_full"
$ws.Range("C4").Value = "resnet18-pretrained, 512-128 linear, relu, 128-5 linear, softmax"
$ws.Range("D4").Value = "pretrained  model weights, all layers trained"
$ws.Range("E4").Value = "Adam"
$ws.Range("F4").Value = "NLLLoss"
$ws.Range("G4").Value = "RandomResizedCrop(300)"
$ws.Range("H4").Value = "Dataset mean, std normalise"
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 0.0004
$ws.Range("K4").Value = "One cycle lr"
$ws.Range("L4").Value = "StratifiedKFold (5 splits)"
$ws.Range("M4").Value = "[0.69351, 0.59728, 0.45618]"
$ws.Range("N4").Value = "[0.75349, 0.82322, 0.84584]"
$ws.Range("O4").Value = "[0.67943, 0.52374, 0.44371]"

# --- Leave the selection where data entry would naturally end up ---
$ws.Range("O5").Select()
